$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 376 - this shifts the existing rows 376..474 down to 377..475,
# matching the dimension change from A1:R474 to A1:R475.
$ws.Rows.Item(376).Insert()

# Populate the newly inserted row 376 with the new record.
$ws.Cells.Item(376, 1).Value = 3
$ws.Cells.Item(376, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(376, 3).Value = "Coquimbo"
$ws.Cells.Item(376, 4).Value = 44841
$ws.Cells.Item(376, 5).Value = 5
$ws.Cells.Item(376, 6).Value = 100112017
$ws.Cells.Item(376, 7).Value = "Apio"
$ws.Cells.Item(376, 8).Value = "Americana (o)"
$ws.Cells.Item(376, 9).Value = "Primera"
$ws.Cells.Item(376, 10).Value = 230
$ws.Cells.Item(376, 11).Value = 8500
$ws.Cells.Item(376, 12).Value = 9000
$ws.Cells.Item(376, 13).Value = 8761
$ws.Cells.Item(376, 14).Value = "`$/docena de matas"
$ws.Cells.Item(376, 15).Value = "Pan de Az$([char]0xFA)car"
$ws.Cells.Item(376, 16).Value = 1460
$ws.Cells.Item(376, 17).Value = 6
$ws.Cells.Item(376, 18).Value = "Hortaliza"
